# New weekly price report: insert a new "Ajo" (garlic) record as the most
# recent observation for Terminal La Palmera de La Serena. The new record
# is inserted at row 215 (just after the existing 2021-04-14 row), pushing
# every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 215; Excel shifts rows 215..278 down to 216..279 and
# carries formatting down from the row above (so D215 keeps the date style).
$ws.Rows.Item(215).EntireRow.Insert()

$ws.Cells.Item(215, 1).Value  = 8
$ws.Cells.Item(215, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(215, 3).Value  = "Coquimbo"
$ws.Cells.Item(215, 4).Value  = 44722
$ws.Cells.Item(215, 5).Value  = 4
$ws.Cells.Item(215, 6).Value  = 100112003
$ws.Cells.Item(215, 7).Value  = "Ajo"
$ws.Cells.Item(215, 8).Value  = "Chino"
$ws.Cells.Item(215, 9).Value  = "Primera"
$ws.Cells.Item(215, 10).Value = 400
$ws.Cells.Item(215, 11).Value = 19000
$ws.Cells.Item(215, 12).Value = 20000
$ws.Cells.Item(215, 13).Value = 19500
$ws.Cells.Item(215, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(215, 15).Value = "China"
$ws.Cells.Item(215, 16).Value = 1950
$ws.Cells.Item(215, 17).Value = 10
$ws.Cells.Item(215, 18).Value = "Hortaliza"
